# Auto-generated edit script: update cryptos list (Price / Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '53.914.00'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '2.253.47'
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  -0.16%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '494.05'
$cell.ClearFormats()
$ws.Range("E5").Value = '  -0.42%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '128.14'
$cell.ClearFormats()
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  -0.99%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.0946'
$cell.ClearFormats()
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("E10").Value = '  +0.97%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.334'
$cell.ClearFormats()
$ws.Range("E11").Value = '  +3.06%  '
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("D13").Value = '2.653.07'
$ws.Range("E13").Value = '  -0.70%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '22.65'
$cell.ClearFormats()
$ws.Range("E14").Value = '  +4.04%  '
$ws.Range("D15").Value = '53.873.57'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").Value = '2.250.01'
$ws.Range("E17").Value = '  -1.84%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '10.20'
$cell.ClearFormats()
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("E19").Value = '  +0.21%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '300.13'
$cell.ClearFormats()
$ws.Range("E20").Value = '  -0.71%  '
$ws.Range("E21").Value = '  -3.02%  '
$ws.Range("E22").Value = '  +0.44%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '60.68'
$cell.ClearFormats()
$ws.Range("E23").Value = '  -2.81%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  -1.97%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '7.26'
$cell.ClearFormats()
$ws.Range("E26").Value = '  +2.35%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '170.24'
$cell.ClearFormats()
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("E28").Value = '  -0.56%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '5.92'
$cell.ClearFormats()
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("D30").Value = '0.0₃0683'
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("E32").Value = '  -0.02%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '17.74'
$cell.ClearFormats()
$ws.Range("E33").Value = '  +0.63%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.ClearFormats()
$ws.Range("E34").Value = '  +0.54%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.939'
$cell.ClearFormats()
$ws.Range("E35").Value = '  +8.00%  '
$ws.Range("E36").Value = '  -0.65%  '
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("E39").Value = '  -1.53%  '
$ws.Range("E40").Value = '  +0.00%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '125.38'
$cell.ClearFormats()
$ws.Range("E41").Value = '  -2.41%  '
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("E43").Value = '  +0.78%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0887'
$cell.ClearFormats()
$ws.Range("E44").Value = '  -0.38%  '
$ws.Range("E45").Value = '  -0.91%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '238.83'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("E47").Value = '  -1.40%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0203'
$cell.ClearFormats()
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  +0.37%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '16.09'
$cell.ClearFormats()
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("E51").Value = '  -1.04%  '
